# "Finished support for multiple plate in CFrenaming"
# Adds a second plate's data: a "plate2" label on Sheet3 (the source sample
# sheet) and a corresponding "Plate2P1" label on Sheet1 (the renamed /
# output sheet), both placed in cell B11. Finally Sheet1 becomes the
# active sheet.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("Sheet1")
$ws3 = $wb.Worksheets.Item("Sheet3")

# Sheet3 gets the new "plate2" entry first (it ends up earlier in the
# shared-strings table), and its selection moves onto the new cell.
$ws3.Range("B11").Value = "plate2"
$ws3.Activate()
$ws3.Range("B11").Select()

# Sheet1 gets the matching "Plate2P1" entry for the second plate.
$ws1.Range("B11").Value = "Plate2P1"

# Sheet1 ends up being the active sheet/selection when the workbook is saved.
$ws1.Activate()
$ws1.Range("H23").Select()
